# data.xlsx — "Ecole" row 1 edit:
#   C1: "IUT" -> "U"
#   E1: numeric 7764 -> text "7764" (same displayed value, stored as text)
#   G1: "-0.588054" -> "-0.5880540000000565" (more precise longitude, stored as text)
#
# E1/G1 must become literal text cells (not formulas, not auto-coerced back to
# numbers) while keeping their original cell style. Assigning a numeric-looking
# string straight to .Value / .Formula gets auto-coerced to a Number by Excel,
# so instead we write a formula that evaluates to the exact text we want, then
# Copy + PasteSpecial(xlPasteValues) to freeze it into a literal value without
# disturbing the cell's existing number format / style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C1: plain text replacement
$ws.Range("C1").Value = "U"

# E1: force literal text "7764"
$ws.Range("E1").Formula = '="7764"'
$ws.Range("E1").Copy()
$ws.Range("E1").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

# G1: force literal text "-0.5880540000000565"
$ws.Range("G1").Formula = '="-0.5880540000000565"'
$ws.Range("G1").Copy()
$ws.Range("G1").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false
